$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 0.07975464681371225
$ws.Range("C3").Value = -3.017636378339217
$ws.Range("E3").Value = -3.229247082222797
$ws.Range("C4").Value = 0.003352386816724007
$ws.Range("E4").Value = -1.461031976610316
$ws.Range("C5").Value = 3.925837669383347
$ws.Range("E5").Value = 2.641604203902781
$ws.Range("C6").Value = 0.2381541440396262
$ws.Range("E6").Value = 1.60268309892857
$ws.Range("C7").Value = 4.993892964711621
$ws.Range("E7").Value = 2.260118192030736
$ws.Range("C8").Value = 6.711795724673664
$ws.Range("E8").Value = 6.409878804372982
$ws.Range("C9").Value = 0.5121603413743347
$ws.Range("E9").Value = 3.290935868252554
$ws.Range("C10").Value = 1.745565778643887
$ws.Range("E10").Value = 0.7985845180024986
$ws.Range("C11").Value = 2.687500891103922
$ws.Range("E11").Value = 1.922191950024699
$ws.Range("C12").Value = 3.654655474034474
$ws.Range("E12").Value = 3.068403604789749
$ws.Range("C13").Value = 3.712036718632117
$ws.Range("E13").Value = 3.908921577463587
$ws.Range("C14").Value = 2.849400388885992
$ws.Range("E14").Value = 3.535456592693387
$ws.Range("C15").Value = -4.741003096464214
$ws.Range("E15").Value = -2.156795995006056
$ws.Range("C16").Value = 1.194925448553708
$ws.Range("E16").Value = -2.068675356622807
$ws.Range("C17").Value = -1.084365158506884
$ws.Range("E17").Value = -1.089896342664354
$ws.Range("C18").Value = -3.40787540386569
$ws.Range("E18").Value = -2.191935020614488
$ws.Range("C19").Value = -1.853660925652212
$ws.Range("E19").Value = -1.878672029998096
